$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "select max(retailprice) from product;"

$ws.Range("A18").Value = " select id from product where retailprice_cs = 12.5"

$dataRow18 = $ws.Range("B18:I18")
$dataRow18.NumberFormat = "@"
$ws.Range("B18").Value = "0.00068527"
$ws.Range("C18").Value = "0.00048519"
$ws.Range("D18").Value = "0.00048550"
$ws.Range("E18").Value = "0.00106290"
$ws.Range("F18").Value = "0.00049251"
$ws.Range("G18").Value = "0.00047968"
$ws.Range("H18").Value = "0.00048119"
$ws.Range("I18").Value = "0.00053873"
$dataRow18.ClearFormats()

$ws.Range("A27").Select()
